$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $ws.Range("J2").Value = "0.5,0.0625"
    $ws.Range("J2").Select()
}

$wb.Worksheets.Item("NATURE").Activate()
$wb.Worksheets.Item("NATURE").Range("J2").Select()
